# Sheet1 holds an IMPACT_code / usda_code / IMPACT_conversion lookup table.
# The row for IMPACT_code "ctool" (usda_code "04670") is being removed -
# it was row 44; rows below it (cwhea, cyams) shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Locate the row whose column A holds "ctool" so this isn't hard-coded
# to a row number that might drift.
$targetRow = 0
$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row  # xlDown
for ($r = 1; $r -le 46; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "ctool") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows($targetRow).Delete()
}

# Re-apply the sort on the (now one-row-shorter) table so the worksheet's
# recorded sort range shrinks along with the data (A2:C46 -> A2:C45).
$lastDataRow = 45
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A" + $lastDataRow))
$sortObj.SetRange($ws.Range("A2:C" + $lastDataRow))
$sortObj.Header = 0
$sortObj.Apply()

# Reflect the row-delete selection state: the freshly-deleted row's
# position (now occupied by the following row) ends up selected as a
# whole row, with the view scrolled near the top of the table.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A44:XFD44").Select()
